# Applies updated Benders decomposition results (fixed mistake in code)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -2593924273.244156
$ws.Range("D3").Value = 12.0644725412605
$ws.Range("F3").Value = 4.437355799993384

$ws.Range("C4").Value = -2593924273.244156
$ws.Range("D4").Value = 12.0644725412605
$ws.Range("F4").Value = 4.064896199997747

$ws.Range("C5").Value = -2624353816.835186
$ws.Range("D5").Value = 18.09670881189074
$ws.Range("F5").Value = 4.27215960000467

$ws.Range("C6").Value = -2624353816.835186
$ws.Range("D6").Value = 18.09670881189074
$ws.Range("F6").Value = 4.232226200008881

$ws.Range("C7").Value = -3141656057.882689
$ws.Range("D7").Value = 120.644725412605
$ws.Range("F7").Value = 4.536767899990082

$ws.Range("C8").Value = -3141656057.882689
$ws.Range("D8").Value = 120.644725412605
$ws.Range("F8").Value = 4.135162399994442

$ws.Range("C9").Value = -3750246929.70328
$ws.Range("D9").Value = 241.2894508252099
$ws.Range("F9").Value = 3.980953700011014

$ws.Range("C10").Value = -3750246929.703279
$ws.Range("D10").Value = 241.2894508252099
$ws.Range("F10").Value = 4.252947700006189

$ws.Range("C11").Value = -3794187190.648727
$ws.Range("D11").Value = 250
$ws.Range("F11").Value = 3.977119200004381

$ws.Range("C12").Value = -3794187190.648727
$ws.Range("D12").Value = 250
$ws.Range("F12").Value = 4.197983900012332

$ws.Range("C13").Value = -3794187190.648727
$ws.Range("D13").Value = 250
$ws.Range("F13").Value = 3.975355799993849

$ws.Range("C14").Value = -3794187190.648727
$ws.Range("D14").Value = 250
$ws.Range("F14").Value = 4.294374000004609

$ws.Range("C15").Value = -3794187190.648727
$ws.Range("D15").Value = 250
$ws.Range("F15").Value = 4.375703600002453

$ws.Range("C16").Value = -3794187190.648727
$ws.Range("D16").Value = 250
$ws.Range("F16").Value = 4.057205500008422

$ws.Range("C17").Value = -3794187190.648727
$ws.Range("D17").Value = 250
$ws.Range("F17").Value = 4.173780500001158

$ws.Range("C18").Value = -3794187190.648727
$ws.Range("D18").Value = 250
$ws.Range("F18").Value = 4.184106500004418

$ws.Range("C19").Value = -3794187190.648727
$ws.Range("D19").Value = 250
$ws.Range("F19").Value = 4.062346799997613

$ws.Range("C20").Value = -3794187190.648727
$ws.Range("D20").Value = 250
$ws.Range("F20").Value = 4.317703799999435

$ws.Range("F21").Value = 4.023150499997428

$ws.Range("F22").Value = 4.262797599993064
